# Update the "Fitness" values (column C) on the active worksheet to reflect
# the new run results (rows 2-152, i.e. generations 0-150).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C2").Value = 11878
$ws.Range("C3:C5").Value = 11399
$ws.Range("C6:C7").Value = 11153
$ws.Range("C8:C8").Value = 10246
$ws.Range("C9:C10").Value = 9957
$ws.Range("C11:C11").Value = 9653
$ws.Range("C12:C13").Value = 9301
$ws.Range("C14:C19").Value = 8829
$ws.Range("C20:C27").Value = 8531
$ws.Range("C28:C29").Value = 8510
$ws.Range("C30:C42").Value = 8235
$ws.Range("C43:C48").Value = 8202
$ws.Range("C49:C64").Value = 8101
$ws.Range("C65:C91").Value = 8068
$ws.Range("C92:C152").Value = 7573
